$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text first so numeric-looking price strings (e.g. trailing
# zeros, thousand-dot formats) are preserved exactly as authored, not coerced
# into numbers by Excel's automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "25.976.20"
$ws.Range("E2").Value = "  -1.59%  "

# Row 3
$ws.Range("D3").Value = "1.637.62"
$ws.Range("E3").Value = "  -1.82%  "

# Row 4
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.57%  "

# Row 5
$ws.Range("D5").Value = "215.90"
$ws.Range("E5").Value = "  -1.45%  "

# Row 6
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "1.014"
$ws.Range("E6").Value = "  +0.68%  "

# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.5003"
$ws.Range("E7").Value = "  -3.42%  "

# Row 8
$ws.Range("D8").Value = "0.2572"
$ws.Range("E8").Value = "  -0.24%  "

# Row 9
$ws.Range("D9").Value = "0.06425"
$ws.Range("E9").Value = "  -0.59%  "

# Row 10
$ws.Range("D10").Value = "19.48"
$ws.Range("E10").Value = "  -2.54%  "

# Row 11
$ws.Range("D11").Value = "0.07762"
$ws.Range("E11").Value = "  +1.18%  "

# Row 12
$ws.Range("D12").Value = "1.649.57"
$ws.Range("E12").Value = "  -1.35%  "

# Row 13
$ws.Range("D13").Value = "4.256"
$ws.Range("E13").Value = "  -2.09%  "

# Row 14
$ws.Range("D14").Value = "1.865.92"
$ws.Range("E14").Value = "  -1.65%  "

# Row 15
$ws.Range("D15").Value = "0.5451"
$ws.Range("E15").Value = "  -1.99%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7937"
$ws.Range("E16").Value = "  -1.46%  "

# Row 17
$ws.Range("D17").Value = "63.60"
$ws.Range("E17").Value = "  -1.79%  "

# Row 18
$ws.Range("D18").Value = "26.000.50"
$ws.Range("E18").Value = "  -1.67%  "

# Row 19
$ws.Range("D19").Value = "1.014"
$ws.Range("E19").Value = "  +0.65%  "

# Row 20
$ws.Range("D20").Value = "203.76"
$ws.Range("E20").Value = "  -3.05%  "

# Row 21
$ws.Range("D21").Value = "4.308"
$ws.Range("E21").Value = "  -2.66%  "

# Row 22
$ws.Range("D22").Value = "10.01"
$ws.Range("E22").Value = "  -1.13%  "

# Row 23
$ws.Range("D23").Value = "5.958"
$ws.Range("E23").Value = "  +1.02%  "

# Row 24
$ws.Range("D24").Value = "1.015"
$ws.Range("E24").Value = "  +0.70%  "

# Row 25
$ws.Range("E25").Value = "  +14.00%  "

# Row 26
$ws.Range("D26").Value = "141.28"
$ws.Range("E26").Value = "  -2.62%  "

# Row 27
$ws.Range("D27").Value = "0.1150"
$ws.Range("E27").Value = "  -1.61%  "

# Row 28
$ws.Range("D28").Value = "15.75"
$ws.Range("E28").Value = "  -0.40%  "

# Row 29
$ws.Range("D29").Value = "6.803"
$ws.Range("E29").Value = "  -2.96%  "

# Row 30
$ws.Range("D30").Value = "0.05047"
$ws.Range("E30").Value = "  -3.95%  "

# Row 31
$ws.Range("D31").Value = "1.244"
$ws.Range("E31").Value = "  -1.56%  "

# Row 32
$ws.Range("D32").Value = "3.267"
$ws.Range("E32").Value = "  -3.43%  "

# Row 33
$ws.Range("D33").Value = "3.202"
$ws.Range("E33").Value = "  -0.71%  "

# Row 34
$ws.Range("D34").Value = "1.546"
$ws.Range("E34").Value = "  -2.04%  "

# Row 35
$ws.Range("D35").Value = "2.352"
$ws.Range("E35").Value = "  -1.07%  "

# Row 36
$ws.Range("D36").Value = "0.8926"
$ws.Range("E36").Value = "  -4.00%  "

# Row 37
$ws.Range("D37").Value = "2.620"
$ws.Range("E37").Value = "  -5.08%  "

# Row 38
$ws.Range("D38").Value = "0.5648"
$ws.Range("E38").Value = "  -1.33%  "

# Row 39
$ws.Range("D39").Value = "1.133.79"
$ws.Range("E39").Value = "  -1.57%  "

# Row 40
$ws.Range("D40").Value = "0.01562"
$ws.Range("E40").Value = "  -2.85%  "

# Row 41
$ws.Range("D41").Value = "2.582"
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("D42").Value = "1.014"
$ws.Range("E42").Value = "  +0.71%  "

# Row 43
$ws.Range("D43").Value = "5.652"
$ws.Range("E43").Value = "  -0.08%  "

# Row 44
$ws.Range("D44").Value = "0.8192"
$ws.Range("E44").Value = "  -3.06%  "

# Row 45
$ws.Range("D45").Value = "99.85"
$ws.Range("E45").Value = "  -0.44%  "

# Row 46
$ws.Range("D46").Value = "1.776.95"
$ws.Range("E46").Value = "  -1.61%  "

# Row 47
$ws.Range("E47").Value = "  +2.99%  "

# Row 48
$ws.Range("D48").Value = "0.4541"
$ws.Range("E48").Value = "  +1.03%  "

# Row 49
$ws.Range("D49").Value = "1.014"
$ws.Range("E49").Value = "  +0.66%  "

# Row 50
$ws.Range("D50").Value = "54.84"
$ws.Range("E50").Value = "  -2.17%  "

# Row 51
$ws.Range("D51").Value = "0.05030"
$ws.Range("E51").Value = "  -1.65%  "

# Restore the default cell style on column D so no stray number-format / style
# index is left behind now that the text values are safely stored.
$ws.Range("D2:D51").Style = "Normal"

Write-Host "done"